# Updates cryptos list data per upstream refresh (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.270.39'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.667.06'
$ws.Range("E3").Value = '  +3.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.90'
$ws.Range("E5").Value = '  +4.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.49'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.668.01'
$ws.Range("E9").Value = '  +3.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.66'
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.357'
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.37'
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.141.31'
$ws.Range("E15").Value = '  +3.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.110.48'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000145'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.663.50'
$ws.Range("E18").Value = '  +3.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.46'
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.82'
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.39'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.86'
$ws.Range("E22").Value = '  +3.26%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.63'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("E26").Value = '  -1.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.165'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.53'
$ws.Range("E28").Value = '  +3.09%  '
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '542.22'
$ws.Range("E29").Value = '  +18.17%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.87'
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  +5.11%  '
$ws.Range("E33").Value = '  +8.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0810'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '172.38'
$ws.Range("E35").Value = '  -2.32%  '
$ws.Range("E36").Value = '  +13.09%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.24'
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.83'
$ws.Range("E40").Value = '  +7.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.13'
$ws.Range("E41").Value = '  +9.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.74'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.05'
$ws.Range("E44").Value = '  +3.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0564'
$ws.Range("E45").Value = '  +4.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.634'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0963'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0240'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.80'
$ws.Range("E49").Value = '  +4.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.73'
$ws.Range("E50").Value = '  +2.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.32'
$ws.Range("E51").Value = '  -0.77%  '
